$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")
$ws.Range("B2").Value = "TestNew"
$ws.Range("B2").Select()
